$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a genuine TEXT value into a cell without disturbing that
# cell's existing number format / style. We stage the text in a scratch
# cell that is explicitly formatted as Text, copy it, and Paste Special
# "Values" into the destination -- which carries over the string type but
# leaves the destination's own formatting (and the workbook's style table)
# untouched.
# ---------------------------------------------------------------------------
$scratch = $ws.Cells.Item(500, 20)

function Set-TextValue {
    param($targetCell, [string]$text)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)  # xlPasteValues
}

# Rows 183-189 currently carry style index 4 (no-number-format / right
# top alignment); the edit upgrades A/B on those rows to style index 5,
# matching rows 176-182 and the rest of the table. Pull that formatting
# across with Paste Special "Formats" so no new style entries are created.
$fmtSrc = $ws.Range("A175:B175")
$fmtSrc.Copy()
$ws.Range("A183:B189").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# New order rows appended to the tracking table.
# Columns: A = Remessa, B = Material, C = Quantidade
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 176; A = "80266955"; B = "10359-ARI-I"; C = 8 },
    @{ Row = 177; A = "80266955"; B = "10381-ARI-I"; C = 4 },
    @{ Row = 178; A = "80266955"; B = "10253-ARI-I"; C = 1 },
    @{ Row = 179; A = "80266955"; B = "10035-ARI-I"; C = 1 },
    @{ Row = 180; A = "80266955"; B = "10257-ARI-I"; C = 4 },
    @{ Row = 181; A = "80266955"; B = "10636-ARI-I"; C = 2 },
    @{ Row = 182; A = "80266955"; B = "10645-ARI-I"; C = 4 },
    @{ Row = 183; A = "80266955"; B = "10637-ARI-I"; C = 4 },
    @{ Row = 184; A = "80266959"; B = "10078-BLB-I"; C = 20 },
    @{ Row = 185; A = "80266968"; B = "33640-ATE-I"; C = 6 },
    @{ Row = 186; A = "80266968"; B = "33378-ATE-I"; C = 5 },
    @{ Row = 187; A = "80266968"; B = "33380-ATE-I"; C = 5 },
    @{ Row = 188; A = "80266981"; B = "10499-ARI-I"; C = 3 },
    @{ Row = 189; A = "80266981"; B = "10381-ARI-I"; C = 1 }
)

foreach ($entry in $rows) {
    Set-TextValue $ws.Cells.Item($entry.Row, 1) $entry.A
    Set-TextValue $ws.Cells.Item($entry.Row, 2) $entry.B
    $ws.Cells.Item($entry.Row, 3).Value = $entry.C
}

$scratch.Clear()
$excel.CutCopyMode = 0

$ws.Range("G10").Select()
